# Insert 3 new rows at position 738, pushing the existing rows 738-835 down to 741-838.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("738:740").Insert()

function Set-DataRow {
    param(
        $Row,
        $Fecha,
        $Variedad,
        $Calidad,
        $Volumen,
        $PrecioMin,
        $PrecioMax,
        $PrecioProm,
        $Unidad,
        $Origen,
        $PrecioKg,
        $KgUnidades
    )

    $ws.Cells.Item($Row, 1).Value2  = 10
    $ws.Cells.Item($Row, 2).Value2  = "Vega Modelo de Temuco"
    $ws.Cells.Item($Row, 3).Value2  = "La Araucanía"
    $ws.Cells.Item($Row, 4).Value2  = $Fecha
    $ws.Cells.Item($Row, 5).Value2  = 9
    $ws.Cells.Item($Row, 6).Value2  = 100112043
    $ws.Cells.Item($Row, 7).Value2  = "Pepino ensalada"
    $ws.Cells.Item($Row, 8).Value2  = $Variedad
    $ws.Cells.Item($Row, 9).Value2  = $Calidad
    $ws.Cells.Item($Row, 10).Value2 = $Volumen
    $ws.Cells.Item($Row, 11).Value2 = $PrecioMin
    $ws.Cells.Item($Row, 12).Value2 = $PrecioMax
    $ws.Cells.Item($Row, 13).Value2 = $PrecioProm
    $ws.Cells.Item($Row, 14).Value2 = $Unidad
    $ws.Cells.Item($Row, 15).Value2 = $Origen
    $ws.Cells.Item($Row, 16).Value2 = $PrecioKg
    $ws.Cells.Item($Row, 17).Value2 = $KgUnidades
    $ws.Cells.Item($Row, 18).Value2 = "Hortaliza"
}

Set-DataRow 738 45127 "Alaska" "Primera" 100 28000 28000 28000 "$/caja 50 unidades" "Región de Arica y Parinacota" 560 50

Set-DataRow 739 45127 "Sin especificar" "Primera" 500 14000 15000 14600 "$/caja 50 unidades" "Región de Arica y Parinacota" 292 50

Set-DataRow 740 45127 "Sin especificar" "Segunda" 180 10000 12000 11111 "$/caja 60 unidades" "Región de Arica y Parinacota" 185 60
